$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (existing row) - rename test id and refresh wording
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "OPQA-1102"
$ws.Range("B2").Value = "Verify that system is able to recommend peoples for user"
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = "1PRECOMMEND"
$ws.Range("D2").Value = "/recommend/people/(SYS_USER1)"
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("E2").Value = "GET"
$ws.Range("G2").Value = "?max=6"
$ws.Range("J2").Value = "status=200"

# ---------------------------------------------------------------------------
# Row 3 (existing row)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "OPQA-1103"
$ws.Range("B3").Value = "Verify that system is able to recommend articles for user"
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = "1PRECOMMEND"
$ws.Range("D3").Value = "/recommend/articles/(SYS_USER1)"
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("E3").Value = "GET"
$ws.Range("G3").Value = "?max=3"
$ws.Range("J3").Value = "status=200"

# ---------------------------------------------------------------------------
# Row 4 (new row)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "OPQA-1399"
$ws.Range("B4").Value = "Verify that 1P-recommend API endpoint for article recommendations based on user selected documents"
$ws.Range("B4").WrapText = $true
$ws.Range("C4").Value = "1PRECOMMEND"
$ws.Range("D4").Value = "/recommend/debug/articles/(SYS_USER1)"
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("E4").Value = "GET"
$ws.Range("G4").Value = "?max=3"
$ws.Range("J4").Value = "status=200"
$ws.Rows.Item(4).RowHeight = 45

# ---------------------------------------------------------------------------
# Row 5 (new row)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "OPQA-1400"
$ws.Range("B5").Value = "Verify that 1P-recommend API endpoint to obtain total times cited for a given ORCID"
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = "1PRECOMMEND"
$ws.Range("D5").Value = "/recommend/jcrmetrix=rid&query=0000-0002-1553-596x"
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("E5").Value = "GET"
$ws.Range("J5").Value = "status=200"
$ws.Rows.Item(5).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 6 (new row) - 1PSEARCHV3 search test
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "OPQA-896"
$ws.Range("A6").WrapText = $true
$ws.Range("B6").Value = "Verify that to get articles for query"
$ws.Range("B6").WrapText = $true
$ws.Range("C6").Value = "1PSEARCHV3"
$ws.Range("D6").Value = "/wos/search"
$ws.Range("D6").Style = "Hyperlink"
$ws.Range("E6").Value = "GET"
$ws.Range("G6").Value = "?query=biotechnology&size=1&fields=citingsrcslocalcount&sort=citingsrcslocalcount:desc"
$ws.Range("G6").WrapText = $true
$ws.Range("J6").Value = "status=200"
$ws.Range("J6").WrapText = $true
$ws.Range("K6").Value = "hits.hits._id"
$ws.Range("K6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 45

# ---------------------------------------------------------------------------
# Row 7 (new row)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "OPQA-1401"
$ws.Range("B7").Value = "Verify that user should receive article recommendation on an article page"
$ws.Range("B7").WrapText = $true
$ws.Range("C7").Value = "1PRECOMMEND"
$ws.Range("D7").Value = "/recommend/matchingdocs/(OPQA-896_hits.hits._id)"
$ws.Range("D7").Style = "Hyperlink"
$ws.Range("E7").Value = "GET"
$ws.Range("G7").Value = "?source=articles&fields=title"
$ws.Range("J7").Value = "status=200"
$ws.Rows.Item(7).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 8 (new row)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "OPQA-1402"
$ws.Range("B8").Value = "Verify that 1P-recommend API endpoint for predicted categories"
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value = "1PRECOMMEND"
$ws.Range("D8").Value = "/recommend/predict/biotechnology"
$ws.Range("D8").Style = "Hyperlink"
$ws.Range("E8").Value = "GET"
$ws.Range("J8").Value = "status=200"
$ws.Rows.Item(8).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 9 (new row)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "OPQA-1403"
$ws.Range("B9").Value = "Verify that 1P-recommend API endpoint for recommending articles and peoples "
$ws.Range("B9").WrapText = $true
$ws.Range("C9").Value = "1PRECOMMEND"
$ws.Range("D9").Value = "/recommend/(SYS_USER1)"
$ws.Range("D9").Style = "Hyperlink"
$ws.Range("E9").Value = "GET"
$ws.Range("J9").Value = "status=200"
$ws.Rows.Item(9).RowHeight = 30

# ---------------------------------------------------------------------------
# Sheet level tweaks
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.57
$ws.Range("B9").Select()

Write-Output "done"
